# Adiciona o sobrenome "Valença" logo apos "Ivna" no paragrafo do cliente,
# criando um novo run com a mesma formatacao (Arial - fonte complexa,
# negrito, italico, tamanho 40 (20pt) e idioma pt-BR).

$d = $word.ActiveDocument

# Localiza o texto "Ivna" dentro do documento.
$rng = $d.Content
$found = $rng.Find.Execute("Ivna", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Texto 'Ivna' nao encontrado no documento."
}

# Colapsa o range encontrado para o seu final (logo apos "Ivna").
$rng.Collapse(0)

# Insere o sobrenome como um novo run logo depois de "Ivna".
$rng.InsertAfter("Valença")

# Aplica ao novo run a mesma formatacao usada em "Ivna":
# fonte complexa Arial, negrito, italico, tamanho 40 (20pt) e pt-BR.
$rng.Font.NameBi = "Arial"
$rng.Font.Bold = $true
$rng.Font.BoldBi = $true
$rng.Font.Italic = $true
$rng.Font.ItalicBi = $true
$rng.Font.Size = 20
$rng.Font.SizeBi = 20
$rng.LanguageID = "pt-BR"
